$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 3.142102271962365
$ws.Range("D2").Value = 3.712465629373006
$ws.Range("E2").Value = 16.58665469366118
$ws.Range("F2").Value = 19.16919346126199
$ws.Range("G2").Value = 20.25933481174998
$ws.Range("H2").Value = 11.91216814505267
$ws.Range("I2").Value = 16.28600980451986
$ws.Range("K2").Value = 12.23378974099224
$ws.Range("N2").Value = 15.84976004048714
$ws.Range("O2").Value = 17.01519987945291
$ws.Range("C3").Value = 3.026146567211337
$ws.Range("D3").Value = 3.643005623447176
$ws.Range("E3").Value = 15.64065884266869
$ws.Range("F3").Value = 19.12824881328886
$ws.Range("G3").Value = 20.13930893388924
$ws.Range("H3").Value = 11.95096719524141
$ws.Range("I3").Value = 16.29720475172894
$ws.Range("K3").Value = 11.58480487203048
$ws.Range("N3").Value = 15.8605722277719
$ws.Range("O3").Value = 17.05059954428774
$ws.Range("C4").Value = 2.953148265520523
$ws.Range("D4").Value = 3.59916126006985
$ws.Range("E4").Value = 15.03436956984698
$ws.Range("F4").Value = 19.10995798842047
$ws.Range("G4").Value = 20.0752927672116
$ws.Range("H4").Value = 11.97728645026426
$ws.Range("I4").Value = 16.30911125797749
$ws.Range("K4").Value = 11.16505001818093
$ws.Range("N4").Value = 15.86888461912918
$ws.Range("O4").Value = 17.07752834424415
$ws.Range("C5").Value = 2.92299593851342
$ws.Range("D5").Value = 3.581007073503454
$ws.Range("E5").Value = 14.78117064480914
$ws.Range("F5").Value = 19.10423032858943
$ws.Range("G5").Value = 20.05166239513616
$ws.Range("H5").Value = 11.98863806179244
$ws.Range("I5").Value = 16.31522558043885
$ws.Range("K5").Value = 10.98874386231791
$ws.Range("N5").Value = 15.87269404036762
$ws.Range("O5").Value = 17.08980223897905
$ws.Range("C6").Value = 2.917966206439717
$ws.Range("D6").Value = 3.577975639333516
$ws.Range("E6").Value = 14.73876570188793
$ws.Range("F6").Value = 19.10338356485412
$ws.Range("G6").Value = 20.04788751905344
$ws.Range("H6").Value = 11.99056077191811
$ws.Range("I6").Value = 16.31631699064336
$ws.Range("K6").Value = 10.95915436157769
$ws.Range("N6").Value = 15.87335211716854
$ws.Range("O6").Value = 17.09191866121908
$ws.Range("C7").Value = 2.952743193324593
$ws.Range("D7").Value = 3.598917571280808
$ws.Range("E7").Value = 15.03097927354105
$ws.Range("F7").Value = 19.10987375153528
$ws.Range("G7").Value = 20.0749641069698
$ws.Range("H7").Value = 11.97743700817646
$ws.Range("I7").Value = 16.30918861164222
$ws.Range("K7").Value = 11.16269341596407
$ws.Range("N7").Value = 15.86893428389117
$ws.Range("O7").Value = 17.07768861782773
$ws.Range("C8").Value = 3.102521847453984
$ws.Range("D8").Value = 3.688771812350798
$ws.Range("E8").Value = 16.26590072440844
$ws.Range("F8").Value = 19.15365703960962
$ws.Range("G8").Value = 20.21595603309101
$ws.Range("H8").Value = 11.92502704685101
$ws.Range("I8").Value = 16.28882358866415
$ws.Range("K8").Value = 12.01448775236387
$ws.Range("N8").Value = 15.85314144102567
$ws.Range("O8").Value = 17.0263245539008
$ws.Range("C9").Value = 3.424334111125984
$ws.Range("D9").Value = 3.854934309424067
$ws.Range("E9").Value = 18.59110249843191
$ws.Range("F9").Value = 19.29359732231236
$ws.Range("G9").Value = 20.56798145624823
$ws.Range("H9").Value = 11.84212374052569
$ws.Range("I9").Value = 16.28892739528953
$ws.Range("K9").Value = 13.51271075690047
$ws.Range("N9").Value = 15.83539872530913
$ws.Range("O9").Value = 16.96703707160918
$ws.Range("C10").Value = 3.6512504824965
$ws.Range("D10").Value = 3.970199243640252
$ws.Range("E10").Value = 20.24224550517165
$ws.Range("F10").Value = 19.42887624275601
$ws.Range("G10").Value = 20.87059419312544
$ws.Range("H10").Value = 11.79341513703977
$ws.Range("I10").Value = 16.31350491454826
$ws.Range("K10").Value = 14.5051302588383
$ws.Range("N10").Value = 15.83035548700587
$ws.Range("O10").Value = 16.94901986553418
$ws.Range("C11").Value = 3.749277725985803
$ws.Range("D11").Value = 4.021023831025534
$ws.Range("E11").Value = 20.95086522886039
$ws.Range("F11").Value = 19.49731369484669
$ws.Range("G11").Value = 21.01726815435886
$ws.Range("H11").Value = 11.77392304993485
$ws.Range("I11").Value = 16.33000699279138
$ws.Range("K11").Value = 14.9326172765868
$ws.Range("N11").Value = 15.82978071425017
$ws.Range("O11").Value = 16.94641648638225
$ws.Range("C12").Value = 3.785641055007503
$ws.Range("D12").Value = 4.040028085047146
$ws.Range("E12").Value = 21.2131269118154
$ws.Range("F12").Value = 19.52420488566492
$ws.Range("G12").Value = 21.07405443175547
$ws.Range("H12").Value = 11.76692659765829
$ws.Range("I12").Value = 16.33701962281538
$ws.Range("K12").Value = 15.0910212281048
$ws.Range("N12").Value = 15.82980886637718
$ws.Range("O12").Value = 16.94623763249715
$ws.Range("C13").Value = 3.777843387264358
$ws.Range("D13").Value = 4.035946096081074
$ws.Range("E13").Value = 21.1569139246415
$ws.Range("F13").Value = 19.5183702901628
$ws.Range("G13").Value = 21.06177002974261
$ws.Range("H13").Value = 11.76841627319342
$ws.Range("I13").Value = 16.33547539943541
$ws.Range("K13").Value = 15.05706109095966
$ws.Range("N13").Value = 15.82979189454685
$ws.Range("O13").Value = 16.946240224348
$ws.Range("C14").Value = 3.75228458956751
$ws.Range("D14").Value = 4.02259222778296
$ws.Range("E14").Value = 20.97256326220726
$ws.Range("F14").Value = 19.49950660268377
$ws.Range("G14").Value = 21.02191537453033
$ws.Range("H14").Value = 11.77333972794992
$ws.Range("I14").Value = 16.33056863391228
$ws.Range("K14").Value = 14.94571907254061
$ws.Range("N14").Value = 15.82977811247317
$ws.Range("O14").Value = 16.94638558539098
$ws.Range("C15").Value = 3.736530182786112
$ws.Range("D15").Value = 4.014380797096726
$ws.Range("E15").Value = 20.85885270800023
$ws.Range("F15").Value = 19.48807855055828
$ws.Range("G15").Value = 20.99766364893468
$ws.Range("H15").Value = 11.77640563777
$ws.Range("I15").Value = 16.32766248461812
$ws.Range("K15").Value = 14.87706543898578
$ws.Range("N15").Value = 15.82980163804563
$ws.Range("O15").Value = 16.94657978697374
$ws.Range("C16").Value = 3.644738671942302
$ws.Range("D16").Value = 3.966844434781619
$ws.Range("E16").Value = 20.19508331251818
$ws.Range("F16").Value = 19.42454113818232
$ws.Range("G16").Value = 20.86118571130636
$ws.Range("H16").Value = 11.79474275703974
$ws.Range("I16").Value = 16.3125334355022
$ws.Range("K16").Value = 14.47670796604445
$ws.Range("N16").Value = 15.83042752655027
$ws.Range("O16").Value = 16.94930279463207
$ws.Range("C17").Value = 3.587088421307003
$ws.Range("D17").Value = 3.937262430115609
$ws.Range("E17").Value = 19.7770224397125
$ws.Range("F17").Value = 19.3873187044851
$ws.Range("G17").Value = 20.77973364759756
$ws.Range("H17").Value = 11.80667578788663
$ws.Range("I17").Value = 16.30461442318053
$ws.Range("K17").Value = 14.2249390893638
$ws.Range("N17").Value = 15.83125086180093
$ws.Range("O17").Value = 16.95240789183764
$ws.Range("C18").Value = 3.553440700496469
$ws.Range("D18").Value = 3.920096677766532
$ws.Range("E18").Value = 19.53256321747017
$ws.Range("F18").Value = 19.36655977965954
$ws.Range("G18").Value = 20.73373637585725
$ws.Range("H18").Value = 11.81379017648553
$ws.Range("I18").Value = 16.30056067082999
$ws.Range("K18").Value = 14.07787336842934
$ws.Range("N18").Value = 15.83188635203144
$ws.Range("O18").Value = 16.9547201552269
$ws.Range("C19").Value = 3.541964542721077
$ws.Range("D19").Value = 3.914259048967184
$ws.Range("E19").Value = 19.44910521709793
$ws.Range("F19").Value = 19.35964334831921
$ws.Range("G19").Value = 20.7183103174013
$ws.Range("H19").Value = 11.81624202045021
$ws.Range("I19").Value = 16.29927422941466
$ws.Range("K19").Value = 14.02769317772189
$ws.Range("N19").Value = 15.83212937288906
$ws.Range("O19").Value = 16.9555933438039
$ws.Range("C20").Value = 3.593276058789359
$ws.Range("D20").Value = 3.940427186954106
$ws.Range("E20").Value = 19.82193964882231
$ws.Range("F20").Value = 19.39121389203639
$ws.Range("G20").Value = 20.78831656702964
$ws.Range("H20").Value = 11.80537952621488
$ws.Range("I20").Value = 16.30540556682309
$ws.Range("K20").Value = 14.25197392849536
$ws.Range("N20").Value = 15.83114646755165
$ws.Range("O20").Value = 16.95202285726422
$ws.Range("C21").Value = 3.75981245477859
$ws.Range("D21").Value = 4.026521227410723
$ws.Range("E21").Value = 21.02687619055795
$ws.Range("F21").Value = 19.50502099462588
$ws.Range("G21").Value = 21.03358832322126
$ws.Range("H21").Value = 11.77188313587533
$ws.Range("I21").Value = 16.33198916206314
$ws.Range("K21").Value = 14.9785174591696
$ws.Range("N21").Value = 15.82977550125294
$ws.Range("O21").Value = 16.94632096935995
$ws.Range("C22").Value = 3.864235288404262
$ws.Range("D22").Value = 4.081374276169557
$ws.Range("E22").Value = 21.77896442571121
$ws.Range("F22").Value = 19.58507687425046
$ws.Range("G22").Value = 21.20111221331774
$ws.Range("H22").Value = 11.75223472258116
$ws.Range("I22").Value = 16.35381280680161
$ws.Range("K22").Value = 15.43308569788879
$ws.Range("N22").Value = 15.83031153639665
$ws.Range("O22").Value = 16.94729898728984
$ws.Range("C23").Value = 3.808909922493927
$ws.Range("D23").Value = 4.052230839191911
$ws.Range("E23").Value = 21.38078962638211
$ws.Range("F23").Value = 19.54183614997267
$ws.Range("G23").Value = 21.11105861723301
$ws.Range("H23").Value = 11.76251571283229
$ws.Range("I23").Value = 16.34175874034423
$ws.Range("K23").Value = 15.19233657640191
$ws.Range("N23").Value = 15.82989490258914
$ws.Range("O23").Value = 16.94634580656661
$ws.Range("C24").Value = 3.590480196415584
$ws.Range("D24").Value = 3.938996893781351
$ws.Range("E24").Value = 19.80164538486451
$ws.Range("F24").Value = 19.38945088057273
$ws.Range("G24").Value = 20.7844336378864
$ws.Range("H24").Value = 11.80596477521158
$ws.Range("I24").Value = 16.30504633629613
$ws.Range("K24").Value = 14.23975870315183
$ws.Range("N24").Value = 15.83119315899778
$ws.Range("O24").Value = 16.95219528966388
$ws.Range("C25").Value = 3.336239530784544
$ws.Range("D25").Value = 3.811132233472639
$ws.Range("E25").Value = 17.94515048916015
$ws.Range("F25").Value = 19.24999324768929
$ws.Range("G25").Value = 20.46485634560436
$ws.Range("H25").Value = 11.86241477323953
$ws.Range("I25").Value = 16.2845997320986
$ws.Range("K25").Value = 13.12623099666955
$ws.Range("N25").Value = 15.83878972986224
$ws.Range("O25").Value = 16.97860841979946
